# "Generate Report for Handoff"
# The a2cde370-1b2d-4683-ad17-4171c983b947.md file has moved from
# "In Translation" to "Ready for handoff" for both the zh-cn and de-de
# locales, with fresh handoff timestamps. Reflect that across the
# Overview roll-up sheet and the two per-locale detail sheets, and
# widen the status/locale columns so the longer text still fits
# (mirrors Excel's own column auto-fit after the content grew).

$wb = $excel.ActiveWorkbook

# ---- Overview sheet --------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Row 3 is the a2cde370-...md file. Columns E (zh-cn) / F (de-de) hold
# the per-locale status, column G the latest handoff-xliff-generate date.
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-09-02 18:15:29"

# Columns E & F grew to fit "Ready for handoff" (was "In Translation").
$ws.Columns.Item(5).ColumnWidth = 16.33
$ws.Columns.Item(6).ColumnWidth = 16.33

# ---- zh-cn detail sheet ------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Row 3 is the a2cde370-...md file.
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "mt"
$ws.Range("H3").Value = "2016-09-02 18:15:24"

# Column C (Status) grew to fit the longer status text.
$ws.Columns.Item(3).ColumnWidth = 16.33

# ---- de-de detail sheet ------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

# Row 3 is the a2cde370-...md file.
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "mt"
$ws.Range("H3").Value = "2016-09-02 18:15:29"

# Column C (Status) grew to fit the longer status text.
$ws.Columns.Item(3).ColumnWidth = 16.33
